# Auto-generated script applying cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.257.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "'2.221.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'297.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.65%  "
$ws.Range("D6").Value = "'89.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.42%  "
$ws.Range("D7").Value = "'0.559"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.46%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.490"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.94%  "
$ws.Range("D10").Value = "'32.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.99%  "
$ws.Range("D11").Value = "'0.0777"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.84%  "
$ws.Range("D12").Value = "'6.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.80%  "
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "'2.557.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'2.207.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.44%  "
$ws.Range("D16").Value = "'13.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "'0.776"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.61%  "
$ws.Range("D18").Value = "'43.957.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "0.0₃0905"
$ws.Range("E19").Value = "  -6.27%  "
$ws.Range("D20").Value = "'5.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.11%  "
$ws.Range("D21").Value = "'10.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -10.40%  "
$ws.Range("D22").Value = "'64.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("D23").Value = "'236.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("E24").Value = "  -7.30%  "
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").Value = "'37.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "'9.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.69%  "
$ws.Range("E30").Value = "  -3.42%  "
$ws.Range("D31").Value = "'148.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.67%  "
$ws.Range("D32").Value = "'5.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -10.40%  "
$ws.Range("D33").Value = "'2.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.23%  "
$ws.Range("D34").Value = "'0.0749"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.60%  "
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("E36").Value = "  -6.92%  "
$ws.Range("D37").Value = "'2.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -12.82%  "
$ws.Range("E38").Value = "  -7.41%  "
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("E40").Value = "  -8.15%  "
$ws.Range("E41").Value = "  -8.44%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").Value = "'13.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.35%  "
$ws.Range("D44").Value = "'1.804.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.80%  "
$ws.Range("D45").Value = "'1.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.81%  "
$ws.Range("E46").Value = "  -7.99%  "
$ws.Range("D47").Value = "'74.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.49%  "
$ws.Range("D48").Value = "'93.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.38%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'13.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "'66.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.74%  "
$ws.Range("D51").Value = "'2.441.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.27%  "
